# Update TPM-derived metrics for the Fgf17-Fgfr3 LR-pair sheet.
# Rows 2-4 keep Sending cluster = FAPs (D varies over ECs/FAPs/MuSCs),
# rows 5-7 keep Sending cluster = MuSCs. The ligand-side expression
# stats (E-H) changed only for the MuSCs sending-cluster rows (5-7),
# while downstream derived-specificity / receptor / edge values (I-T)
# were recomputed for every row given the new TPM inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs)
$ws.Range("I2").Value = 0.6594814792829158
$ws.Range("J2").Value = 0.6594814792829158
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.752937333333333
$ws.Range("N2").Value = 11.258812
$ws.Range("O2").Value = 0.6855621274031838
$ws.Range("P2").Value = 0.6855621274031838
$ws.Range("Q2").Value = 0.2223665409164444
$ws.Range("R2").Value = 2.001298868248
$ws.Range("S2").Value = 0.4521155259201944
$ws.Range("T2").Value = 0.4521155259201944

# Row 3 (FAPs -> FAPs)
$ws.Range("I3").Value = 0.6594814792829158
$ws.Range("J3").Value = 0.6594814792829158
$ws.Range("O3").Value = 0.2368266084628361
$ws.Range("P3").Value = 0.2368266084628362
$ws.Range("S3").Value = 0.1561827620826271
$ws.Range("T3").Value = 0.1561827620826271

# Row 4 (FAPs -> MuSCs)
$ws.Range("I4").Value = 0.6594814792829158
$ws.Range("J4").Value = 0.6594814792829158
$ws.Range("M4").Value = 0.4248633333333334
$ws.Range("N4").Value = 1.27459
$ws.Range("O4").Value = 0.07761126413398003
$ws.Range("P4").Value = 0.07761126413398005
$ws.Range("Q4").Value = 0.02517371898444444
$ws.Range("R4").Value = 0.22656347086
$ws.Range("S4").Value = 0.05118319128009426
$ws.Range("T4").Value = 0.05118319128009427

# Row 5 (MuSCs -> ECs)
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.030594
$ws.Range("H5").Value = 0.091782
$ws.Range("I5").Value = 0.3405185207170842
$ws.Range("J5").Value = 0.3405185207170842
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.752937333333333
$ws.Range("N5").Value = 11.258812
$ws.Range("O5").Value = 0.6855621274031838
$ws.Range("P5").Value = 0.6855621274031838
$ws.Range("Q5").Value = 0.114817364776
$ws.Range("R5").Value = 1.033356282984
$ws.Range("S5").Value = 0.2334466014829893
$ws.Range("T5").Value = 0.2334466014829893

# Row 6 (MuSCs -> FAPs)
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.030594
$ws.Range("H6").Value = 0.091782
$ws.Range("I6").Value = 0.3405185207170842
$ws.Range("J6").Value = 0.3405185207170842
$ws.Range("O6").Value = 0.2368266084628361
$ws.Range("P6").Value = 0.2368266084628362
$ws.Range("Q6").Value = 0.039663519914
$ws.Range("R6").Value = 0.356971679226
$ws.Range("S6").Value = 0.08064384638020904
$ws.Range("T6").Value = 0.08064384638020906

# Row 7 (MuSCs -> MuSCs)
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.030594
$ws.Range("H7").Value = 0.091782
$ws.Range("I7").Value = 0.3405185207170842
$ws.Range("J7").Value = 0.3405185207170842
$ws.Range("M7").Value = 0.4248633333333334
$ws.Range("N7").Value = 1.27459
$ws.Range("O7").Value = 0.07761126413398003
$ws.Range("P7").Value = 0.07761126413398005
$ws.Range("Q7").Value = 0.01299826882
$ws.Range("R7").Value = 0.11698441938
$ws.Range("S7").Value = 0.02642807285388577
$ws.Range("T7").Value = 0.02642807285388577
